# Weekly update: insert a new price record for "Vega Monumental Concepción -
# Berenjena" as the latest week, pushing all existing rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 164; this shifts rows 164-195 down to 165-196 and
# extends the sheet dimension to A1:R196, carrying over the date-number
# style (s="2") that column D uses.
$ws.Rows.Item(164).Insert()

# Populate the newly inserted row 164 with the new weekly record.
$ws.Cells.Item(164, 1).Value  = 11
$ws.Cells.Item(164, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(164, 3).Value  = "Bíobío"
$ws.Cells.Item(164, 4).Value  = 45218
$ws.Cells.Item(164, 5).Value  = 8
$ws.Cells.Item(164, 6).Value  = 100112001
$ws.Cells.Item(164, 7).Value  = "Berenjena"
$ws.Cells.Item(164, 8).Value  = "Sin especificar"
$ws.Cells.Item(164, 9).Value  = "Primera"
$ws.Cells.Item(164, 10).Value = 100
$ws.Cells.Item(164, 11).Value = 7500
$ws.Cells.Item(164, 12).Value = 8000
$ws.Cells.Item(164, 13).Value = 7750
$ws.Cells.Item(164, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(164, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(164, 16).Value = 155
$ws.Cells.Item(164, 17).Value = 50
$ws.Cells.Item(164, 18).Value = "Hortaliza"
